# "found problems with dataset" - fix the CRONACA/CRONACA NERA/POLITICA per-social
# stats table on Sheet1: the dataset had a row missing (a bad "CROANCA" row
# inserted at the top of the CRONACA block) which shifted every following
# row down by one. Re-write rows 3-12 with the corrected data and push the
# table's dimension out to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,  "CROANCA",      "Facebook",  3,  33,   3,   1),
    @(4,  "CRONACA",      "Facebook",  11, 1058, 81,  2),
    @(5,  "CRONACA",      "Instagram", 17, 939,  174, 7),
    @(6,  "CRONACA",      "YouTube",   8,  1093, 77,  3),
    @(7,  "CRONACA NERA", "Facebook",  24, 1014, 150, 12),
    @(8,  "CRONACA NERA", "Instagram", 21, 996,  168, 15),
    @(9,  "CRONACA NERA", "YouTube",   12, 1117, 66,  5),
    @(10, "POLITICA",     "Facebook",  14, 1111, 69,  5),
    @(11, "POLITICA",     "Instagram", 38, 1001, 145, 9),
    @(12, "POLITICA",     "YouTube",   21, 1104, 68,  2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("A12").Value = ""
